$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("X20").Value = 476
$ws.Range("AB20").Value = 1739

$ws.Range("X21").Value = 566
$ws.Range("AB21").Value = 2021

$ws.Range("X22").Value = 796
$ws.Range("AB22").Value = 2742

$ws.Range("X23").Value = 1039
$ws.Range("AB23").Value = 3601

$ws.Range("X61").Value = 5203
$ws.Range("AB61").Value = 28602

$ws.Range("X62").Value = 5216
$ws.Range("AB62").Value = 28741

$ws.Range("X63").Value = 5228
$ws.Range("AB63").Value = 28822

$ws.Range("X64").Value = 5252
$ws.Range("AB64").Value = 28987

$ws.Range("X65").Value = 5271
$ws.Range("AB65").Value = 29095

$ws.Range("X66").Value = 5274
$ws.Range("AB66").Value = 29193

$ws.Range("AB67").Value = 29253

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("X59").Value = 340
$ws.Range("AB59").Value = 1577

$ws.Range("X60").Value = 350
$ws.Range("AB60").Value = 1607

$ws.Range("X61").Value = 353
$ws.Range("AB61").Value = 1636

$ws.Range("X62").Value = 363
$ws.Range("AB62").Value = 1660

$ws.Range("X63").Value = 366
$ws.Range("AB63").Value = 1677

$ws.Range("X64").Value = 370
$ws.Range("AB64").Value = 1698

$ws.Range("X65").Value = 375
$ws.Range("AB65").Value = 1723

$ws.Range("X66").Value = 375
$ws.Range("AB66").Value = 1727

$ws.Range("AB67").Value = 1733

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("X27").Value = 209
$ws.Range("AB27").Value = 903

$ws.Range("X28").Value = 229
$ws.Range("AB28").Value = 1089

$ws.Range("X29").Value = 247
$ws.Range("AB29").Value = 1206

$ws.Range("X30").Value = 283
$ws.Range("AB30").Value = 1363

$ws.Range("X31").Value = 311
$ws.Range("AB31").Value = 1462

$ws.Range("X32").Value = 322
$ws.Range("AB32").Value = 1608

$ws.Range("X33").Value = 356
$ws.Range("AB33").Value = 1789

$ws.Range("X34").Value = 364
$ws.Range("AB34").Value = 1876

$ws.Range("X35").Value = 367
$ws.Range("AB35").Value = 1991

$ws.Range("X36").Value = 370
$ws.Range("AB36").Value = 2176

$ws.Range("X37").Value = 378
$ws.Range("AB37").Value = 2204

$ws.Range("X38").Value = 383
$ws.Range("AB38").Value = 2283

$ws.Range("X39").Value = 374
$ws.Range("AB39").Value = 2347

$ws.Range("X40").Value = 365
$ws.Range("AB40").Value = 2328

$ws.Range("X41").Value = 365
$ws.Range("AB41").Value = 2316

$ws.Range("X42").Value = 376
$ws.Range("AB42").Value = 2301

$ws.Range("X43").Value = 361
$ws.Range("AB43").Value = 2306

$ws.Range("X44").Value = 334
$ws.Range("AB44").Value = 2224

$ws.Range("X45").Value = 325
$ws.Range("AB45").Value = 2134

$ws.Range("X46").Value = 308
$ws.Range("AB46").Value = 2066

$ws.Range("X47").Value = 299
$ws.Range("AB47").Value = 2006

$ws.Range("X48").Value = 285
$ws.Range("AB48").Value = 1931

$ws.Range("X49").Value = 294
$ws.Range("AB49").Value = 1908

$ws.Range("X50").Value = 294
$ws.Range("AB50").Value = 1891

$ws.Range("X51").Value = 277
$ws.Range("AB51").Value = 1851

$ws.Range("X52").Value = 251
$ws.Range("AB52").Value = 1729

$ws.Range("X53").Value = 244
$ws.Range("AB53").Value = 1673

$ws.Range("X54").Value = 228
$ws.Range("AB54").Value = 1574

$ws.Range("X55").Value = 225
$ws.Range("AB55").Value = 1523

$ws.Range("X56").Value = 227
$ws.Range("AB56").Value = 1516

$ws.Range("X57").Value = 210
$ws.Range("AB57").Value = 1492

$ws.Range("X58").Value = 197
$ws.Range("AB58").Value = 1414

$ws.Range("X59").Value = 183
$ws.Range("AB59").Value = 1348

$ws.Range("X60").Value = 183
$ws.Range("AB60").Value = 1295

$ws.Range("X61").Value = 171
$ws.Range("AB61").Value = 1249

$ws.Range("X62").Value = 176
$ws.Range("AB62").Value = 1220

$ws.Range("X63").Value = 184
$ws.Range("AB63").Value = 1199

$ws.Range("X64").Value = 175
$ws.Range("AB64").Value = 1177

$ws.Range("X65").Value = 164
$ws.Range("AB65").Value = 1148

$ws.Range("X66").Value = 156
$ws.Range("AB66").Value = 1081

$ws.Range("AB67").Value = 1063

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("X34").Value = 60
$ws.Range("AB34").Value = 280

$ws.Range("X35").Value = 61
$ws.Range("AB35").Value = 290

$ws.Range("X36").Value = 63
$ws.Range("AB36").Value = 336

$ws.Range("X37").Value = 66
$ws.Range("AB37").Value = 358

$ws.Range("X38").Value = 68
$ws.Range("AB38").Value = 384

$ws.Range("X39").Value = 69
$ws.Range("AB39").Value = 390

$ws.Range("X40").Value = 67
$ws.Range("AB40").Value = 396

$ws.Range("X41").Value = 67
$ws.Range("AB41").Value = 396

$ws.Range("X42").Value = 66
$ws.Range("AB42").Value = 386

$ws.Range("X43").Value = 64
$ws.Range("AB43").Value = 386

$ws.Range("X44").Value = 63
$ws.Range("AB44").Value = 378

$ws.Range("X45").Value = 62
$ws.Range("AB45").Value = 379

$ws.Range("X46").Value = 62
$ws.Range("AB46").Value = 372

$ws.Range("X47").Value = 64
$ws.Range("AB47").Value = 365

$ws.Range("X48").Value = 64
$ws.Range("AB48").Value = 365

$ws.Range("X49").Value = 63
$ws.Range("AB49").Value = 355

$ws.Range("X50").Value = 59
$ws.Range("AB50").Value = 350

$ws.Range("X51").Value = 58
$ws.Range("AB51").Value = 335

$ws.Range("X52").Value = 55
$ws.Range("AB52").Value = 315

$ws.Range("X53").Value = 54
$ws.Range("AB53").Value = 294

$ws.Range("X54").Value = 53
$ws.Range("AB54").Value = 288

$ws.Range("X55").Value = 50
$ws.Range("AB55").Value = 275

$ws.Range("X56").Value = 50
$ws.Range("AB56").Value = 268

$ws.Range("X57").Value = 44
$ws.Range("AB57").Value = 254

$ws.Range("X58").Value = 44
$ws.Range("AB58").Value = 242

$ws.Range("X59").Value = 42
$ws.Range("AB59").Value = 207

$ws.Range("X60").Value = 40
$ws.Range("AB60").Value = 199

$ws.Range("X61").Value = 36
$ws.Range("AB61").Value = 196

$ws.Range("X62").Value = 36
$ws.Range("AB62").Value = 191

$ws.Range("X63").Value = 36
$ws.Range("AB63").Value = 181

$ws.Range("X64").Value = 34
$ws.Range("AB64").Value = 178

$ws.Range("X65").Value = 33
$ws.Range("AB65").Value = 164

$ws.Range("X66").Value = 33
$ws.Range("AB66").Value = 170

$ws.Range("AB67").Value = 163
